# Daily "cryptos list" refresh (GitHub Actions bot) - updates the Price (D)
# and Volume(1h) (E) columns for the crypto rows on Sheet1.
#
# Price-column values that look numeric (e.g. "212.27", "0.06300") are
# written with a leading apostrophe text-qualifier and the style is reset
# to "Normal" afterwards so Excel stores them as plain text (matching the
# source data, which keeps trailing zeros / decimal-grouped numbers like
# "26.188.17" as text) without leaving the cell tagged with a different
# style than the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'26.188.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = "'1.676.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  -0.78%  '
$ws.Range("D5").Value = "'212.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.43%  '
$ws.Range("D6").Value = "'0.5261"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.04%  '
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("E8").Value = '  -3.24%  '
$ws.Range("D9").Value = "'0.06300"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -3.28%  '
$ws.Range("D11").Value = "'0.07560"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = "'1.688.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").Value = "'0.5620"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").Value = "'66.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = "'0.000008031"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.57%  '
$ws.Range("D17").Value = "'26.260.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").Value = "'4.822"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("D20").Value = "'187.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("E21").Value = '  -5.09%  '
$ws.Range("D22").Value = "'6.185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D24").Value = "'149.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = "'0.1249"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.02%  '
$ws.Range("D26").Value = "'7.571"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.64%  '
$ws.Range("D27").Value = "'16.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("D28").Value = "'0.06184"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("D29").Value = "'1.360"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").Value = "'1.288"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("E31").Value = '  -3.57%  '
$ws.Range("D32").Value = "'3.440"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.74%  '
$ws.Range("D33").Value = "'1.634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").Value = "'1.002"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("D35").Value = "'0.6069"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.19%  '
$ws.Range("D36").Value = "'2.406"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").Value = "'6.098"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").Value = "'0.01614"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("D40").Value = "'1.083.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.61%  '
$ws.Range("D41").Value = "'0.8724"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("E42").Value = '  -1.15%  '
$ws.Range("D43").Value = "'100.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'1.827.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("D46").Value = "'56.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.03%  '
$ws.Range("D47").Value = "'0.9977"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").Value = "'8.051"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("D49").Value = "'0.05236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("D50").Value = "'0.4256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("D51").Value = "'5.975"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '
